# Turn the GitHub repo URL text on the "GitHub Link" slide into a clickable
# hyperlink pointing at the same address (mirrors selecting the run and
# using Insert > Link > Insert Link in the PowerPoint UI).

$p = $ppt.ActivePresentation

# Slide 11 ("GitHub Link") -> shape 2 ("Content Placeholder 2") holds the
# repo URL text.
$slide = $p.Slides.Item(11)
$shape = $slide.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

$url = "https://github.com/NarendraShende/Stegno_Project.git"

# ppMouseClick = 1
$hyperlink = $textRange.ActionSettings(1).Hyperlink
$hyperlink.Address = $url
